$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (Consulting Data Engineer @ The Cigna Group):
#  - extend end date
#  - recolor css RGB
#  - flip text color to white
$ws.Range("D2").Value = 45139
$ws.Range("G2").Value = "#006688"
$ws.Range("I2").Value = "white"

# Row 4 (Consulting Data Engineer @ Angel Studios): recolor css RGB
$ws.Range("G4").Value = "#4400cd"

# Row 7 (Senior Data Engineer): fix employer name
$ws.Range("B7").Value = "b"

# Row heights recalculated by Excel as a side effect of the above edits
$ws.Rows.Item(7).RowHeight = 212.25
$ws.Rows.Item(8).RowHeight = 88.5
$ws.Rows.Item(9).RowHeight = 212.25
$ws.Rows.Item(10).RowHeight = 184.5
